$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename test case name and update count of related test cases
$ws.Range("A2").Value = "Step1Fields"
$ws.Range("B2").Value = 19

# Add note about the one untested/pending test case
$ws.Range("E2").Value = "One test case is untested and waiting on TBH143"

# Update the selected cell to reflect the edit location
$ws.Range("E3").Select()
